# Refresh the cryptos list figures (Price/Volume columns) to match the
# latest scrape, per the "Updated cryptos list ... GitHub Actions" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some of the new Price values parse as plain numbers to Excels normal
# type inference (e.g. "5.15", "27.00"), but the sheet stores every Price
# / Volume cell as literal text (so things like trailing zeroes in "27.00"
# survive). Force those specific cells to the Text number format first so
# assigning the string does not silently convert them to numbers.
$textCells = @("D5", "D6", "D12", "D13", "D14", "D16", "D19", "D20", "D21", "D26", "D27", "D31", "D33", "D38", "D39", "D43", "D46", "D47", "D48")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Price (column D) and Volume(1h) (column E) updates, row by row.
$ws.Range("D2").Value = "67.280.90"
$ws.Range("E2").Value = "  +0.56%  "

$ws.Range("D3").Value = "2.552.10"
$ws.Range("E3").Value = "  -2.33%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "590.97"
$ws.Range("E5").Value = "  +0.75%  "

$ws.Range("D6").Value = "173.84"
$ws.Range("E6").Value = "  +5.08%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("E8").Value = "  +0.55%  "

$ws.Range("D9").Value = "2.551.11"
$ws.Range("E9").Value = "  -2.28%  "

$ws.Range("E10").Value = "  -0.93%  "

$ws.Range("E11").Value = "  +1.78%  "

$ws.Range("D12").Value = "5.15"
$ws.Range("E12").Value = "  -0.77%  "

$ws.Range("D13").Value = "0.348"
$ws.Range("E13").Value = "  -4.89%  "

$ws.Range("D14").Value = "27.00"
$ws.Range("E14").Value = "  -0.81%  "

$ws.Range("D15").Value = "3.011.64"
$ws.Range("E15").Value = "  -2.66%  "

$ws.Range("D16").Value = "0.0000177"
$ws.Range("E16").Value = "  -1.00%  "

$ws.Range("D17").Value = "67.107.68"
$ws.Range("E17").Value = "  +0.42%  "

$ws.Range("D18").Value = "2.549.91"
$ws.Range("E18").Value = "  -2.81%  "

$ws.Range("D19").Value = "8.04"
$ws.Range("E19").Value = "  +3.35%  "

$ws.Range("D20").Value = "11.34"
$ws.Range("E20").Value = "  -2.75%  "

$ws.Range("D21").Value = "356.32"
$ws.Range("E21").Value = "  +0.45%  "

$ws.Range("E22").Value = "  -1.38%  "

$ws.Range("E23").Value = "  +0.72%  "

$ws.Range("E24").Value = "  +6.08%  "

$ws.Range("E25").Value = "  -0.03%  "

$ws.Range("D26").Value = "69.99"
$ws.Range("E26").Value = "  +0.96%  "

$ws.Range("D27").Value = "10.06"
$ws.Range("E27").Value = "  -4.17%  "

$ws.Range("D28").Value = "2.673.05"
$ws.Range("E28").Value = "  -2.66%  "

$ws.Range("E29").Value = "  -0.02%  "

$ws.Range("D30").Value = "0.0₃0988"
$ws.Range("E30").Value = "  -0.31%  "

$ws.Range("D31").Value = "533.63"
$ws.Range("E31").Value = "  -1.19%  "

$ws.Range("E32").Value = "  +0.21%  "

$ws.Range("D33").Value = "1.35"
$ws.Range("E33").Value = "  +1.51%  "

$ws.Range("E34").Value = "  -0.65%  "

$ws.Range("E35").Value = "  -0.52%  "

$ws.Range("E36").Value = "  +0.08%  "

$ws.Range("E37").Value = "  +0.23%  "

$ws.Range("D38").Value = "156.65"
$ws.Range("E38").Value = "  -1.06%  "

$ws.Range("D39").Value = "18.72"
$ws.Range("E39").Value = "  -0.86%  "

$ws.Range("E40").Value = "  +1.18%  "

$ws.Range("E41").Value = "  -1.75%  "

$ws.Range("E42").Value = "  +0.41%  "

$ws.Range("D43").Value = "5.16"
$ws.Range("E43").Value = "  +0.74%  "

$ws.Range("E44").Value = "  +5.01%  "

$ws.Range("D46").Value = "39.71"
$ws.Range("E46").Value = "  -1.24%  "

$ws.Range("D47").Value = "150.29"

$ws.Range("D48").Value = "0.562"
$ws.Range("E48").Value = "  -2.22%  "

$ws.Range("D49").Value = "0.0₆0278"
$ws.Range("E49").Value = "  -5.26%  "

$ws.Range("E50").Value = "  -0.94%  "

$ws.Range("E51").Value = "  +0.16%  "
